$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 751 is the last existing data row (A751:I751); its formatting is the
# template we stamp onto every newly appended row so the new cells pick up
# the same styles (s="3" dates, s="1" text/numbers, s="2" empty "G").
$template = $ws.Range("A751:I751")
# Row 750 has an empty "Localisation douleur" (G) cell - used as the format
# template for new rows whose G column is also blank.
$emptyG = $ws.Range("G750")

function Add-WellnessRow {
    param(
        [int]$Row,
        [int]$Date,
        [string]$Name,
        [int]$Volume,
        [int]$Intensite,
        [int]$Fatigue,
        [int]$Douleur,
        $Localisation,
        [int]$Plaisir
    )

    $dst = $ws.Range("A" + $Row + ":I" + $Row)
    $template.Copy()
    $dst.PasteSpecial(-4122)

    $ws.Range("A$Row").Value = $Date
    $ws.Range("B$Row").Value = $Name
    $ws.Range("C$Row").Value = $Volume
    $ws.Range("D$Row").Value = $Intensite
    $ws.Range("E$Row").Value = $Fatigue
    $ws.Range("F$Row").Value = $Douleur

    if ($Localisation -ne $null) {
        $ws.Range("G$Row").Value = $Localisation
    } else {
        $emptyG.Copy()
        $ws.Range("G$Row").PasteSpecial(-4122)
    }

    $ws.Range("H$Row").Value = $Plaisir
    $ws.Range("I$Row").Formula = "=C" + $Row + "*D" + $Row
}

Add-WellnessRow 752 46043 "Amir Etien"      73 5 7 3 "Ischio"     7
Add-WellnessRow 753 46043 "Yoann Martelat"  73 4 5 4 "Genou"      5
Add-WellnessRow 754 46043 "Romain Thunet"   73 6 7 6 "Poignet"    3
Add-WellnessRow 755 46043 "Yoan Zouma"      73 4 3 2 "Ischio"     6
Add-WellnessRow 756 46043 "Naim Ighbane"    73 6 7 3 "Courbature" 6
Add-WellnessRow 757 46043 "Kamal Bafounta"  73 6 7 4 "Genou"      7
Add-WellnessRow 758 46043 "Naim Dhib"       73 5 4 8 "Hanche"     4
Add-WellnessRow 759 46043 "Karahali Souaré" 73 6 6 6 "Cheville"   4
Add-WellnessRow 760 46043 "Theo Owono"      73 5 4 0 $null        7
Add-WellnessRow 761 46043 "Hedi Nasri"      73 6 5 3 "Hanche"     6

# Restore the selection the author ended up with after entering the new rows.
$null = $ws.Range("K756").Select()
$excel.ActiveWindow.ScrollRow = 736
$excel.ActiveWindow.ScrollColumn = 1
